$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New year column M: headers / data for the 2023 column (mirrors L).
# ---------------------------------------------------------------------------
$ws.Range("M3").Value2 = 2023
$ws.Range("L3").Copy()
$ws.Range("M3").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("M5").Value2 = 14065.6
$ws.Range("L5").Copy()
$ws.Range("M5").PasteSpecial(-4122)

$ws.Range("M6").Value2 = 7161.9
$ws.Range("L6").Copy()
$ws.Range("M6").PasteSpecial(-4122)

$ws.Range("M7").Value2 = 46.213456944602434
$ws.Range("L7").Copy()
$ws.Range("M7").PasteSpecial(-4122)

$ws.Range("M8").Value2 = 4.4790126265498803E-2
$ws.Range("L8").Copy()
$ws.Range("M8").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Row 4 ("Generation of hazardous waste per person") becomes bold
#    9pt Times New Roman, grouped by its pre-existing formatting blocks so
#    each block folds onto a single shared style, same as the source file.
# ---------------------------------------------------------------------------
$r1 = $ws.Range("A4:C4")
$r1.Font.Name = "Times New Roman"
$r1.Font.Size = 9
$r1.Font.Bold = $true

$r2 = $ws.Range("D4:I4")
$r2.Font.Name = "Times New Roman"
$r2.Font.Size = 9
$r2.Font.Bold = $true

$r3 = $ws.Range("J4")
$r3.Font.Name = "Times New Roman"
$r3.Font.Size = 9
$r3.Font.Bold = $true

$r4 = $ws.Range("K4")
$r4.Font.Name = "Times New Roman"
$r4.Font.Size = 9
$r4.Font.Bold = $true

$r5 = $ws.Range("L4")
$r5.Font.Name = "Times New Roman"
$r5.Font.Size = 9
$r5.Font.Bold = $true

# M4 picks up L4's brand-new bold style, then gets its own value.
$ws.Range("L4").Copy()
$ws.Range("M4").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("M4").Value2 = 1963.9481143272037

# ---------------------------------------------------------------------------
# 3. Row 7 label wraps onto two lines now that the row is taller.
# ---------------------------------------------------------------------------
$ws.Range("A7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 24

# ---------------------------------------------------------------------------
# 4. Row 1 grew taller to fit the header text; columns A:C share one width.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 57
$ws.Range("A1:C1").ColumnWidth = 38.1

# ---------------------------------------------------------------------------
# 5. Drop the stale cell selection left over from the previous save.
# ---------------------------------------------------------------------------
$ws.Range("A1").Select()
